{"js": "// Update the date heading and the division problems in the practice\n// table (2024-01-08 Monday -> 2024-01-09 Tuesday; cell values changed\n// per the new day's worksheet).\n\n// 1) Update the date paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst datePara = paragraphs.items[0];\ndatePara.getRange().insertText(\"2024-01-09 Tuesday\", \"Replace\");\n\n// 2) Update the division-problem table. The table has 20 rows x 5\n// columns, but only every 4th row (0, 4, 8, 12, 16) actually holds a\n// problem; the rest are blank spacer rows. Target each populated row\n// by its row index and set new cell text, preserving existing\n// run/paragraph formatting (cell.value only replaces the text).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValuesByRow = {\n  0: [\"32\u00f75=\", \"55\u00f78=\", \"42\u00f77=\", \"79\u00f72=\", \"82\u00f78=\"],\n  4: [\"66\u00f79=\", \"11\u00f74=\", \"15\u00f73=\", \"12\u00f76=\", \"12\u00f79=\"],\n  8: [\"45\u00f79=\", \"31\u00f72=\", \"54\u00f74=\", \"77\u00f77=\", \"70\u00f77=\"],\n  12: [\"32\u00f73=\", \"29\u00f78=\", \"83\u00f73=\", \"65\u00f78=\", \"37\u00f73=\"],\n  16: [\"76\u00f76=\", \"13\u00f73=\", \"32\u00f77=\", \"86\u00f78=\", \"16\u00f73=\"],\n};\n\nfor (const rowIndex of Object.keys(newValuesByRow)) {\n  const rowValues = newValuesByRow[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    table.getCell(Number(rowIndex), col).value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the division problems in the practice\n# table (2024-01-08 Monday -> 2024-01-09 Tuesday; cell values changed\n# per the new day's worksheet).\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph (first paragraph in the body).\n$d.Paragraphs(1).Range.Text = \"2024-01-09 Tuesday\"\n\n# 2) Update the division-problem table. The table has 20 rows x 5\n# columns, but only every 4th row (1, 5, 9, 13, 17 in Word's 1-based\n# indexing) actually holds a problem; the rest are blank spacer rows.\n# Target each populated row by its row index and set new cell text,\n# preserving existing run/paragraph formatting (setting Range.Text\n# only replaces the text content of the cell).\n$t = $d.Tables(1)\n\n$newValuesByRow = @{\n    1  = @(\"32\u00f75=\", \"55\u00f78=\", \"42\u00f77=\", \"79\u00f72=\", \"82\u00f78=\")\n    5  = @(\"66\u00f79=\", \"11\u00f74=\", \"15\u00f73=\", \"12\u00f76=\", \"12\u00f79=\")\n    9  = @(\"45\u00f79=\", \"31\u00f72=\", \"54\u00f74=\", \"77\u00f77=\", \"70\u00f77=\")\n    13 = @(\"32\u00f73=\", \"29\u00f78=\", \"83\u00f73=\", \"65\u00f78=\", \"37\u00f73=\")\n    17 = @(\"76\u00f76=\", \"13\u00f73=\", \"32\u00f77=\", \"86\u00f78=\", \"16\u00f73=\")\n}\n\nforeach ($rowIndex in $newValuesByRow.Keys) {\n    $rowValues = $newValuesByRow[$rowIndex]\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
